$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.092.34"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.563.42"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.96"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.84"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.54"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.958.03"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "2.557.35"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.05"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "43.120.38"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.13"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.25"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.77"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.19"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.19"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("E31").Value = "  -4.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.97"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.38"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.95"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.49"
$ws.Range("E41").Value = "  -6.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.96"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0304"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("D46").Value = "2.002.09"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").Value = "2.808.78"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "82.64"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.71"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  +1.54%  "
